$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 20:34"

# Country statistics refresh. Rows are sorted descending by column B
# ("Casos totales"); the refreshed counts push three countries past their
# neighbour, so those pairs swap row position as well as getting new data:
#   rows 64/65   Libano overtakes Singapur
#   rows 71/72   Irlanda overtakes Estado de Palestina
#   rows 203/204 Dominica overtakes Guam

# Row 4
$ws.Range("B4").Value = 8183742
$ws.Range("C4").Value = 33699
$ws.Range("D4").Value = 5293880
$ws.Range("E4").Value = 2667533
$ws.Range("G4").Value = 486
$ws.Range("H4").Value = 222329
# Row 5
$ws.Range("B5").Value = 7364994
$ws.Range("C5").Value = 59924
$ws.Range("D5").Value = 6448268
$ws.Range("E5").Value = 804588
$ws.Range("G5").Value = 827
$ws.Range("H5").Value = 112138
# Row 8
$ws.Range("B8").Value = 972958
$ws.Range("C8").Value = 13318
$ws.Range("G8").Value = 140
$ws.Range("H8").Value = 33553
# Row 13
$ws.Range("B13").Value = 809684
$ws.Range("C13").Value = 30621
$ws.Range("D13").Value = 104082
$ws.Range("E13").Value = 672477
$ws.Range("G13").Value = 88
$ws.Range("H13").Value = 33125
# Row 23
$ws.Range("B23").Value = 346391
$ws.Range("C23").Value = 4649
$ws.Range("E23").Value = 54691
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = 9800
# Row 33
$ws.Range("B33").Value = 163650
$ws.Range("C33").Value = 3317
$ws.Range("D33").Value = 136036
$ws.Range("E33").Value = 24842
$ws.Range("G33").Value = 46
$ws.Range("H33").Value = 2772
# Row 64
$ws.Range("A64").Value = "Libano"
$ws.Range("B64").Value = 58745
$ws.Range("C64").Value = 1499
$ws.Range("D64").Value = 25994
$ws.Range("E64").Value = 32250
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 501
# Row 65
$ws.Range("A65").Value = "Singapur"
$ws.Range("B65").Value = 57892
$ws.Range("C65").Value = 3
$ws.Range("D65").Value = 57764
$ws.Range("E65").Value = 100
$ws.Range("H65").Value = 28
# Row 71
$ws.Range("A71").Value = "Irlanda"
$ws.Range("B71").Value = 46429
$ws.Range("C71").Value = 1186
$ws.Range("D71").Value = 23364
$ws.Range("E71").Value = 21227
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 1838
# Row 72
$ws.Range("A72").Value = "Estado de Palestina"
$ws.Range("B72").Value = 46100
$ws.Range("C72").Value = 442
$ws.Range("D72").Value = 39585
$ws.Range("E72").Value = 6114
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 401
# Row 90
$ws.Range("B90").Value = 22170
$ws.Range("C90").Value = 534
$ws.Range("D90").Value = 16758
$ws.Range("E90").Value = 4597
$ws.Range("G90").Value = 6
$ws.Range("H90").Value = 815
# Row 138
$ws.Range("B138").Value = 4285
$ws.Range("C138").Value = 30
$ws.Range("D138").Value = 3910
$ws.Range("E138").Value = 343
# Row 152
$ws.Range("B152").Value = 2807
$ws.Range("C152").Value = 9
$ws.Range("E152").Value = 1462
# Row 189
$ws.Range("B189").Value = 253
$ws.Range("C189").Value = 5
$ws.Range("E189").Value = 34
# Row 203
$ws.Range("A203").Value = "Dominica"
$ws.Range("B203").Value = 33
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 29
$ws.Range("E203").Value = 4
$ws.Range("H203").Value = 0
# Row 204
$ws.Range("A204").Value = "Guam"
$ws.Range("D204").Value = 0
$ws.Range("E204").Value = 31
$ws.Range("H204").Value = 1
